# Implements the "start events group" metrics columns.
# Adds 11 new BPMN metric columns (nCondition ... nStartErrorEventDefinition)
# right after the existing nEventBasedGateway column, with header row styling
# matching the existing header cells and a 0 default count in the data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newHeaders = @(
    "nCondition",
    "nIntermediateThrowEvent",
    "nStartEvent",
    "nStartSignalEventDefinition",
    "nStartConditionalEventDefinition",
    "nStartTimerEventDefinition",
    "nStartMessageEventDefinition",
    "nStartCompensateEventDefinition",
    "nStartCancelEventDefinition",
    "nStartEscalationEventDefinition",
    "nStartErrorEventDefinition"
)

# Existing header columns run from A (1) to Z (26); new columns start at AA (27).
$startCol = 27

# Use the existing header cell (A1) as the formatting template so the new
# header cells pick up the same bold font + border style already used by
# the other header cells (instead of creating a duplicate style entry).
$templateHeaderCell = $ws.Cells.Item(1, 1)
$templateHeaderCell.Copy()

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $col = $startCol + $i
    $headerCell = $ws.Cells.Item(1, $col)
    $headerCell.PasteSpecial(-4122) # xlPasteFormats
    $headerCell.Value = $newHeaders[$i]

    $dataCell = $ws.Cells.Item(2, $col)
    $dataCell.Value = 0
}

$excel.CutCopyMode = 0
